# Fill in the "Bill To" / "Service To" customer details on the invoice.
# Order of writes matters only insofar as it controls the order new shared
# strings are interned in xl/sharedStrings.xml; write Name+Address for each
# party before City/Country so the table comes out as:
#   John Smith, #4 Some place, Jane Smith, Some city, Some country

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bill To: (left column)
$ws.Range("B11").Value = "John Smith"
$ws.Range("B12").Value = "#4 Some place"

# Service To: (right column)
$ws.Range("E11").Value = "Jane Smith"
$ws.Range("E12").Value = "#4 Some place"

# City - same for both parties
$ws.Range("B13").Value = "Some city"
$ws.Range("E13").Value = "Some city"

# Country - same for both parties
$ws.Range("B14").Value = "Some country"
$ws.Range("E14").Value = "Some country"
